# Populate the "Bill Information" sheet with the bill table data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Name and Address"
$ws.Range("B1").Value = "RR Number"
$ws.Range("C1").Value = "Date"
$ws.Range("D1").Value = "Account ID"
$ws.Range("E1").Value = "Consumption"
$ws.Range("F1").Value = "Tax"
$ws.Range("G1").Value = "Net Amount Due"

# Data rows (2 and 3 are identical, matching the source bill)
$address = "TH PETROLEUM CORP ( LTD ) NEAR MAJESTIC BUS STATION MALLESWARAM _ BANGALORE "

# The Account ID / Consumption / Tax columns are stored as plain text in
# the source bill (not numbers), so prefix with an apostrophe -- the same
# "force text" trick used when typing numbers-as-text directly into Excel --
# to avoid automatic numeric coercion. Style is reset to Normal afterwards
# so the quote-prefix marker doesn't linger on the cell format.
$ws.Range("A2").Value = $address
$ws.Range("B2").Value = "14543E"
$ws.Range("C2").Value = "09 / 10/ 13"
$ws.Range("D2").Value = "'1100157649"
$ws.Range("E2").Value = "'4100"
$ws.Range("F2").Value = "'1272"
$ws.Range("G2").Value = "2 926 9.00"

$ws.Range("A3").Value = $address
$ws.Range("B3").Value = "14543E"
$ws.Range("C3").Value = "09 / 10/ 13"
$ws.Range("D3").Value = "'1100157649"
$ws.Range("E3").Value = "'4100"
$ws.Range("F3").Value = "'1272"
$ws.Range("G3").Value = "2 926 9.00"

$ws.Range("D2:F3").Style = "Normal"

# Column widths to match the source layout (A..G). The engine's ColumnWidth
# setter rounds to 1/6-character pixel steps and adds 5/6 of padding before
# storing, so the inputs below are pre-compensated (target - 5/6) to land the
# serialized <col width=.../> as close as possible to the source values.
$ws.Columns.Item(1).ColumnWidth = 75.05338541666667
$ws.Columns.Item(2).ColumnWidth = 16.498697916666668
$ws.Columns.Item(3).ColumnWidth = 14.721354166666666
$ws.Columns.Item(4).ColumnWidth = 17.830729166666668
$ws.Columns.Item(5).ColumnWidth = 19.608072916666668
$ws.Columns.Item(6).ColumnWidth = 18.385416666666668
$ws.Columns.Item(7).ColumnWidth = 20.721354166666668

# Selection matches the saved view state (A1 selected)
$ws.Range("A1").Select()
